$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Adição do campo Ano e migração do campo Data":
# insert a brand-new column before column F (6th column). "Procedência" (E)
# stays where it is; "Data" and every column after it shift one place to the
# right, and the freshly inserted column F gets the new "Ano" header.
$ws.Columns.Item(6).Insert()

# The inserted column picks up the width of the column to its left ("E"),
# same as Excel's own "Insert Column" behaviour.
$ws.Columns.Item(6).ColumnWidth = 10.5

$ws.Range("F1").Value = "Ano"

# Leave the active selection on the newly added header cell.
$ws.Range("F1").Select() | Out-Null
